{"js": "// Apply the \"Added many more features\" edit to the Money Heat review doc.\n// Each replacement targets a unique, exact source string so a simple\n// search + insertText(Replace) round-trip reproduces the diff precisely\n// without disturbing surrounding run/paragraph formatting.\n\nconst replacements = [\n  // Title (Heading1) and the later bold restatement of the same title\n  // both change to the identical new text.\n  {\n    find: \"Play Money Heat Slot Game for Free - Review & Pros/Cons\",\n    replace: \"Play Money Heat Slot Free - Review\",\n  },\n  // \"What we like\" bullet list\n  {\n    find: \"Wide betting range\",\n    replace: \"Wide betting range for risk-takers\",\n  },\n  {\n    find: \"High potential for big rewards\",\n    replace: \"Stacked Wild for potential big wins\",\n  },\n  {\n    find: \"Stacked Wild feature\",\n    replace: \"High volatility for bigger rewards\",\n  },\n  // \"What we don't like\" bullet list\n  {\n    find: \"Lower RTP compared to other online slots\",\n    replace: \"Modest graphics compared to other slots\",\n  },\n  {\n    find: \"Dated graphics\",\n    replace: \"Lower RTP compared to other games\",\n  },\n  // Italicized meta-description line near the end of the document\n  {\n    find:\n      \"Read our review of Money Heat online slot game: pros, cons and gameplay. Play for free and discover its features, bonuses, payout, and volatility.\",\n    replace:\n      \"Read our review of Money Heat slot game. Play for free and discover its features and bonuses.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit to the Money Heat review doc.\n# Each Find/Replace targets a unique, exact source string (MatchCase on,\n# wdReplaceAll) so the run reproduces the diff precisely without touching\n# unrelated lowercase occurrences (e.g. \"wide betting range\" in body text)\n# or any other run/paragraph formatting.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# Title (Heading1) and the later bold restatement of the same title both\n# change to the identical new text - ReplaceAll handles both occurrences.\nReplace-Text \"Play Money Heat Slot Game for Free - Review & Pros/Cons\" \"Play Money Heat Slot Free - Review\"\n\n# \"What we like\" bullet list\nReplace-Text \"Wide betting range\" \"Wide betting range for risk-takers\"\nReplace-Text \"High potential for big rewards\" \"Stacked Wild for potential big wins\"\nReplace-Text \"Stacked Wild feature\" \"High volatility for bigger rewards\"\n\n# \"What we don't like\" bullet list\nReplace-Text \"Lower RTP compared to other online slots\" \"Modest graphics compared to other slots\"\nReplace-Text \"Dated graphics\" \"Lower RTP compared to other games\"\n\n# Italicized meta-description line near the end of the document\nReplace-Text \"Read our review of Money Heat online slot game: pros, cons and gameplay. Play for free and discover its features, bonuses, payout, and volatility.\" \"Read our review of Money Heat slot game. Play for free and discover its features and bonuses.\"\n"}
